$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update existing Scotland row (row 6) with revised figures ---
$ws.Cells.Item(6,3).Value = 518
$ws.Cells.Item(6,4).Value = 302
$ws.Cells.Item(6,5).Value = 184
$ws.Cells.Item(6,6).Value = 16

# --- Add new country rows: Finland (row 18) and United States (row 19) ---

# Seed formatting for the new rows by copying an existing data row's used
# range down (keeps date/hyperlink styles + number formats intact), then
# overwrite with the real values and clear out any columns that don't
# apply to the new country.
$ws.Range("A17:J17").Copy($ws.Range("A18:J18"))
$ws.Range("A17:J17").Copy($ws.Range("A19:J19"))

# Row 18: Finland
$ws.Cells.Item(18,1).Value = "Finland"
$ws.Cells.Item(18,2).Value = 43979
$ws.Cells.Item(18,3).Value = 313
$ws.Cells.Item(18,4).ClearContents()
$ws.Cells.Item(18,5).Value = 141
$ws.Cells.Item(18,6).ClearContents()
$ws.Cells.Item(18,7).ClearContents()
$ws.Cells.Item(18,8).ClearContents()
$ws.Cells.Item(18,9).Formula = "=E18/C18"
$ws.Hyperlinks.Add($ws.Cells.Item(18,10), "https://thl.fi/en/web/infectious-diseases/what-s-new/coronavirus-covid-19-latest-updates/situation-update-on-coronavirus", "Coronavirus-related_deaths", [Type]::Missing, "https://thl.fi/en/web/infectious-diseases/what-s-new/coronavirus-covid-19-latest-updates/situation-update-on-coronavirus - Coronavirus-related_deaths") | Out-Null
# The visible cell text should be the full URL including the in-page
# fragment; this does not disturb the hyperlink's stored display text.
$ws.Cells.Item(18,10).Value = "https://thl.fi/en/web/infectious-diseases/what-s-new/coronavirus-covid-19-latest-updates/situation-update-on-coronavirus#Coronavirus-related_deaths"
# Re-apply the standard hyperlink formatting (Hyperlinks.Add mutates the
# cell's style record when the cell already carried the Hyperlink style).
$ws.Cells.Item(18,10).Style = "Hyperlink"

# Row 19: United States
$ws.Cells.Item(19,1).Value = "United States"
$ws.Cells.Item(19,2).Value = 43974
$ws.Cells.Item(19,3).Value = 81372
$ws.Cells.Item(19,4).Value = 55903
$ws.Cells.Item(19,5).Value = 20083
$ws.Cells.Item(19,6).Value = 4247
$ws.Cells.Item(19,7).ClearContents()
$ws.Cells.Item(19,8).Value = 1110
$ws.Cells.Item(19,9).Formula = "=E19/C19"
$ws.Cells.Item(19,10).Value = "https://www.cdc.gov/nchs/nvss/vsrr/covid_weekly/index.htm"
$ws.Hyperlinks.Add($ws.Cells.Item(19,10), "https://www.cdc.gov/nchs/nvss/vsrr/covid_weekly/index.htm") | Out-Null
$ws.Cells.Item(19,10).Style = "Hyperlink"

$ws.Range("C13").Select()
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("J20").Select()
